# Insert a new data row at row 574 (shifts existing rows 574-614 down to 575-615),
# matching the weekly-update commit that adds one more price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 574, pushing everything else down.
$ws.Rows("574").Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A574").Value = 3
$ws.Range("B574").Value = "Femacal de La Calera"
$ws.Range("C574").Value = "Coquimbo"
$ws.Range("D574").Value = 44931
$ws.Range("E574").Value = 5
$ws.Range("F574").Value = 100112037
$ws.Range("G574").Value = "Cebollín"
$ws.Range("H574").Value = "Sin especificar"
$ws.Range("I574").Value = "Primera"
$ws.Range("J574").Value = 280
$ws.Range("K574").Value = 3500
$ws.Range("L574").Value = 4000
$ws.Range("M574").Value = 3786
$ws.Range("N574").Value = '$/paquete 36 unidades'
$ws.Range("O574").Value = "Provincia de Quillota"
$ws.Range("P574").Value = 105
$ws.Range("Q574").Value = 36
$ws.Range("R574").Value = "Hortaliza"

# Match the date-number format used by the rest of column D.
$ws.Range("D574").NumberFormat = "YYYY-MM-DD HH:MM:SS"
